$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Adam Benson's hours (row 6)
$ws.Range("C6").Value = 88.02

# Update Kyle Brown's hours (row 7)
$ws.Range("C7").Value = 33.8

# Delete the row for Alex Keller (row 16) entirely - shifts everything below up
$ws.Rows("16:16").Delete()
